# "Generate Report for Archive"
# Status moves from "Ready for handoff" -> "In Translation" on the Overview
# sheet (zh-cn / de-de status cells) and on each per-locale sheet's Status
# column, then the now-narrower Status columns are re-sized to fit the
# shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Resize the two status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
